$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the new row is stored as text, not numbers, matching the original data's type.
# Must set the Text number format BEFORE assigning values, so numeric-looking strings
# (like "11", "0", "100.00") stay text instead of being coerced to numbers.
$ws.Range("A3:K3").NumberFormat = "@"

# Add a new row 3 that duplicates row 2's data (per commit: "json file code updated")
$ws.Cells.Item(3, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(3, 2).Value = " October 29 2020"
$ws.Cells.Item(3, 3).Value = "Super Kings won by 6 wickets"
$ws.Cells.Item(3, 4).Value = "Kolkata Knight Riders"
$ws.Cells.Item(3, 5).Value = "Chennai Super Kings"
$ws.Cells.Item(3, 6).Value = "Rinku Singh "
$ws.Cells.Item(3, 7).Value = "11"
$ws.Cells.Item(3, 8).Value = "11"
$ws.Cells.Item(3, 9).Value = "1"
$ws.Cells.Item(3, 10).Value = "0"
$ws.Cells.Item(3, 11).Value = "100.00"
